# Auto-generated Excel COM-interop script applying the Zalera_Profits.xlsx diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 110.333336
$ws.Range("I6").Value = 110.333336
$ws.Range("K6").Value = 331.000008
$ws.Range("M6").Value = -219.000008
$ws.Range("H8").Value = 12495
$ws.Range("I8").Value = 13327
$ws.Range("K8").Value = 39981
$ws.Range("M8").Value = -39842
$ws.Range("H69").Value = 8000
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 30000
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -29126
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 8000
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 90000
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -85632
$ws.Range("N72").Value = -62736
$ws.Range("H88").Value = 972.5
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 972.5
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2711.1714
$ws.Range("I132").Value = 1824.4
$ws.Range("J132").Value = 8031.8
$ws.Range("K132").Value = 5473.200000000001
$ws.Range("L132").Value = 24095.4
$ws.Range("M132").Value = -2943.200000000001
$ws.Range("N132").Value = -29155.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 503381.75
$ws.Range("I86").Value = 4510
$ws.Range("J86").Value = 669672.3
$ws.Range("K86").Value = 4510
$ws.Range("L86").Value = 669672.3
$ws.Range("M86").Value = -3387
$ws.Range("N86").Value = -671918.3
$ws.Range("H89").Value = 503381.75
$ws.Range("I89").Value = 4510
$ws.Range("J89").Value = 669672.3
$ws.Range("K89").Value = 22550
$ws.Range("L89").Value = 3348361.5
$ws.Range("M89").Value = -16934
$ws.Range("N89").Value = -3359593.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 73401.64999999999
$ws.Range("I132").Value = 2021.8667
$ws.Range("K132").Value = 6065.6001
$ws.Range("M132").Value = -3535.6001
$ws.Range("H134").Value = 4170.3223
$ws.Range("I134").Value = 3795.196
$ws.Range("K134").Value = 11385.588
$ws.Range("M134").Value = -8850.588

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19043.777
$ws.Range("I70").Value = 15361.75
$ws.Range("K70").Value = 15361.75
$ws.Range("M70").Value = -15091.75
$ws.Range("H73").Value = 19043.777
$ws.Range("I73").Value = 15361.75
$ws.Range("K73").Value = 15361.75
$ws.Range("M73").Value = -14425.75
$ws.Range("H80").Value = 4827.7
$ws.Range("J80").Value = 5092.25
$ws.Range("L80").Value = 5092.25
$ws.Range("N80").Value = -7088.25
$ws.Range("H83").Value = 4827.7
$ws.Range("J83").Value = 5092.25
$ws.Range("L83").Value = 25461.25
$ws.Range("N83").Value = -35445.25
$ws.Range("H107").Value = 873.5
$ws.Range("I107").Value = 904.5
$ws.Range("J107").Value = 749.5
$ws.Range("K107").Value = 904.5
$ws.Range("L107").Value = 749.5
$ws.Range("M107").Value = 1015.5
$ws.Range("N107").Value = -4589.5
$ws.Range("H122").Value = 7333.694
$ws.Range("I122").Value = 9001.919
$ws.Range("J122").Value = 2190
$ws.Range("K122").Value = 27005.757
$ws.Range("L122").Value = 6570
$ws.Range("M122").Value = -24555.757
$ws.Range("N122").Value = -11470
$ws.Range("H126").Value = 3052.72
$ws.Range("I126").Value = 2060.353
$ws.Range("J126").Value = 5161.5
$ws.Range("K126").Value = 6181.059
$ws.Range("L126").Value = 15484.5
$ws.Range("M126").Value = -3711.059
$ws.Range("N126").Value = -20424.5
$ws.Range("H132").Value = 6054.933
$ws.Range("J132").Value = 11199.2
$ws.Range("L132").Value = 33597.60000000001
$ws.Range("N132").Value = -38657.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 100000000
$ws.Range("J43").Value = 100000000
$ws.Range("L43").Value = 100000000
$ws.Range("N43").Value = -100000386
$ws.Range("H61").Value = 798.6
$ws.Range("I61").Value = 798.6
$ws.Range("K61").Value = 798.6
$ws.Range("M61").Value = -596.6
$ws.Range("H68").Value = 10667.333
$ws.Range("I68").Value = 3250.6667
$ws.Range("J68").Value = 25500.666
$ws.Range("K68").Value = 3250.6667
$ws.Range("L68").Value = 25500.666
$ws.Range("M68").Value = -2501.6667
$ws.Range("N68").Value = -26998.666
$ws.Range("H71").Value = 10667.333
$ws.Range("I71").Value = 3250.6667
$ws.Range("J71").Value = 25500.666
$ws.Range("K71").Value = 16253.3335
$ws.Range("L71").Value = 127503.33
$ws.Range("M71").Value = -12509.3335
$ws.Range("N71").Value = -134991.33
$ws.Range("H100").Value = 19233200
$ws.Range("I100").Value = 62501450
$ws.Range("K100").Value = 62501450
$ws.Range("M100").Value = -62500909
$ws.Range("H113").Value = 798.6
$ws.Range("I113").Value = 798.6
$ws.Range("K113").Value = 798.6
$ws.Range("M113").Value = 1371.4
$ws.Range("H132").Value = 6861
$ws.Range("I132").Value = 2427.8572
$ws.Range("J132").Value = 12033
$ws.Range("K132").Value = 7283.571599999999
$ws.Range("L132").Value = 36099
$ws.Range("M132").Value = -4753.571599999999
$ws.Range("N132").Value = -41159

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 74199
$ws.Range("J21").Value = 74199
$ws.Range("L21").Value = 74199
$ws.Range("N21").Value = -74669
$ws.Range("H35").Value = 74199
$ws.Range("J35").Value = 74199
$ws.Range("N35").Value = -74779
$ws.Range("H81").Value = 13723.235
$ws.Range("I81").Value = 2165.8333
$ws.Range("J81").Value = 16199.821
$ws.Range("K81").Value = 4331.6666
$ws.Range("L81").Value = 32399.642
$ws.Range("M81").Value = -3270.6666
$ws.Range("N81").Value = -34521.642
$ws.Range("H84").Value = 13723.235
$ws.Range("I84").Value = 2165.8333
$ws.Range("J84").Value = 16199.821
$ws.Range("K84").Value = 21658.333
$ws.Range("L84").Value = 161998.21
$ws.Range("M84").Value = -16354.333
$ws.Range("N84").Value = -172606.21
$ws.Range("H122").Value = 3999.5
$ws.Range("I122").Value = 3999.5
$ws.Range("K122").Value = 11998.5
$ws.Range("M122").Value = -9548.5
